$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 '63.303.14'
$ws.Cells.Item(2, 5).Value = '  +6.01%  '
Set-TextValue 3 4 '3.119.18'
$ws.Cells.Item(3, 5).Value = '  +3.60%  '
$ws.Cells.Item(4, 5).Value = '  -0.12%  '
Set-TextValue 5 4 '585.75'
$ws.Cells.Item(5, 5).Value = '  +3.79%  '
Set-TextValue 6 4 '144.86'
$ws.Cells.Item(6, 5).Value = '  +3.10%  '
$ws.Cells.Item(7, 5).Value = '  -0.04%  '
Set-TextValue 8 4 '3.111.54'
$ws.Cells.Item(8, 5).Value = '  +3.81%  '
$ws.Cells.Item(9, 5).Value = '  +1.28%  '
$ws.Cells.Item(10, 5).Value = '  +12.31%  '
Set-TextValue 11 4 '5.81'
$ws.Cells.Item(11, 5).Value = '  +9.64%  '
$ws.Cells.Item(12, 5).Value = '  +2.70%  '
Set-TextValue 13 4 '0.0000250'
$ws.Cells.Item(13, 5).Value = '  +7.54%  '
Set-TextValue 14 4 '35.66'
$ws.Cells.Item(14, 5).Value = '  +4.69%  '
$ws.Cells.Item(15, 5).Value = '  +0.20%  '
Set-TextValue 16 4 '3.632.36'
$ws.Cells.Item(16, 5).Value = '  +3.62%  '
Set-TextValue 17 4 '7.18'
$ws.Cells.Item(17, 5).Value = '  -0.75%  '
Set-TextValue 18 4 '63.186.27'
$ws.Cells.Item(18, 5).Value = '  +5.85%  '
Set-TextValue 19 4 '3.113.13'
$ws.Cells.Item(19, 5).Value = '  +3.57%  '
Set-TextValue 20 4 '467.13'
$ws.Cells.Item(20, 5).Value = '  +6.29%  '
Set-TextValue 21 4 '14.08'
$ws.Cells.Item(21, 5).Value = '  +3.10%  '
$ws.Cells.Item(22, 5).Value = '  +0.93%  '
Set-TextValue 23 4 '7.54'
$ws.Cells.Item(23, 5).Value = '  +5.84%  '
Set-TextValue 24 4 '13.31'
$ws.Cells.Item(24, 5).Value = '  -1.96%  '
$ws.Cells.Item(25, 5).Value = '  +1.66%  '
$ws.Cells.Item(26, 5).Value = '  -0.03%  '
$ws.Cells.Item(27, 5).Value = '  +0.53%  '
$ws.Cells.Item(28, 5).Value = '  +4.88%  '
Set-TextValue 30 4 '8.28'
$ws.Cells.Item(30, 5).Value = '  +5.80%  '
Set-TextValue 31 4 '6.83'
$ws.Cells.Item(31, 5).Value = '  +8.55%  '
Set-TextValue 32 4 '27.04'
$ws.Cells.Item(32, 5).Value = '  +4.32%  '
$ws.Cells.Item(33, 5).Value = '  +4.14%  '
Set-TextValue 34 4 '0.0₃0870'
$ws.Cells.Item(34, 5).Value = '  +10.84%  '
Set-TextValue 35 4 '2.41'
$ws.Cells.Item(35, 5).Value = '  +14.97%  '
$ws.Cells.Item(36, 5).Value = '  +4.34%  '
$ws.Cells.Item(37, 2).Value = 'dogwifhat'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 37 4 '3.32'
$ws.Cells.Item(37, 5).Value = '  +18.77%  '
$ws.Cells.Item(38, 2).Value = 'Filecoin'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 38 4 '6.06'
$ws.Cells.Item(38, 5).Value = '  +1.86%  '
Set-TextValue 39 4 '50.91'
$ws.Cells.Item(39, 5).Value = '  +3.68%  '
Set-TextValue 40 4 '433.21'
$ws.Cells.Item(40, 5).Value = '  +7.17%  '
Set-TextValue 41 4 '8.73'
$ws.Cells.Item(41, 5).Value = '  +1.49%  '
Set-TextValue 42 4 '2.921.58'
$ws.Cells.Item(42, 5).Value = '  +5.37%  '
$ws.Cells.Item(43, 5).Value = '  +3.96%  '
$ws.Cells.Item(44, 5).Value = '  +10.09%  '
$ws.Cells.Item(45, 5).Value = '  +4.93%  '
Set-TextValue 46 4 '2.18'
$ws.Cells.Item(46, 5).Value = '  +5.77%  '
Set-TextValue 47 4 '35.20'
$ws.Cells.Item(47, 5).Value = '  +1.94%  '
Set-TextValue 49 4 '124.03'
$ws.Cells.Item(49, 5).Value = '  +0.45%  '
$ws.Cells.Item(50, 5).Value = '  +0.44%  '
Set-TextValue 51 4 '24.56'
$ws.Cells.Item(51, 5).Value = '  +3.05%  '
